# Update the four tracking-number cells (R2:R5) on the "Input" sheet with
# new shipment/tracking numbers, keeping them as text values (as they were
# originally authored) rather than letting Excel auto-coerce the digit
# strings into numbers.
#
# Typing a plain digit string into a General-formatted cell would normally
# be auto-converted to a number by Excel, which would lose the original
# text typing (t="s") of these cells. To reproduce a genuine text value
# without disturbing the existing cell formatting (border/fill/number
# format), we write a text-producing formula and then collapse it down to
# a literal value via Copy + PasteSpecial(values only) - exactly the
# "paste values" trick users use to turn a formula result into plain text
# in place, leaving the cell's style untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

$updates = @{
    "R2" = "51524551"
    "R3" = "51524552"
    "R4" = "51524553"
    "R5" = "51524554"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.Formula = "=""" + $updates[$addr] + """"
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

$excel.CutCopyMode = 0
